$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data pattern per language: code, method_seq, is_active
# PWD / 1 / TRUE, OTP / 2 / TRUE, FINGERPRINT / 3 / TRUE, IRIS / 4 / TRUE, FACE / 5 / TRUE
$codes = @("PWD", "OTP", "FINGERPRINT", "IRIS", "FACE")

$langs = @("fra", "ara")

$row = 7
foreach ($lang in $langs) {
    for ($i = 0; $i -lt $codes.Length; $i++) {
        $ws.Cells.Item($row, 1).Value = $lang
        $ws.Cells.Item($row, 2).Value = $codes[$i]
        $ws.Cells.Item($row, 3).Value = ($i + 1)
        $ws.Cells.Item($row, 4).Value = "TRUE"
        $row = $row + 1
    }
}

# Copy styles from row 6 (A6:D6) to the newly added rows 7-16
$ws.Range("A6:D6").Copy() | Out-Null
$ws.Range("A7:D16").PasteSpecial(-4122) | Out-Null

# Set column widths for B and C
$ws.Columns.Item(2).ColumnWidth = 20.28
$ws.Columns.Item(3).ColumnWidth = 20.84

# Set the active cell selection to A12
$ws.Range("A12").Select() | Out-Null
